# "Pimped out the design iteration function a little bit"
#
# This edit adds five new parameters to the "Aircraft parameters" sheet
# (MTOW, Nprops, Nblades, DL, A) used by the DesignIter function, tweaks
# the horizontal tail-volume-coefficient-ish B10 value from 1 to 1.4, and
# leaves the workbook with the "Aircraft parameters" tab active/selected
# instead of "Mission".

$wb = $excel.ActiveWorkbook

$wsMission = $wb.Worksheets.Item("Mission")
$wsAircraft = $wb.Worksheets.Item("Aircraft parameters")

# --- Aircraft parameters sheet: tweak existing value -----------------
$wsAircraft.Range("B10").Value = 1.4

# --- Aircraft parameters sheet: append new parameter rows ------------
# Values/labels are written in the same order the original author typed
# them (new shared-string entries must land in that exact sequence), then
# the remaining, already-known strings ("kg", "-") and numbers are filled
# in afterwards.

# Row 16: MTOW
$wsAircraft.Range("A16").Value = "MTOW"
$wsAircraft.Range("D16").Value = "Maximum take-off weight"

# Row 17: Nprops
$wsAircraft.Range("A17").Value = "Nprops"
$wsAircraft.Range("D17").Value = "Amount of VTOL propellers"

# Row 18: Nblades
$wsAircraft.Range("A18").Value = "Nblades"
$wsAircraft.Range("D18").Value = "Amount of blades on the propeller"

# Row 19: DL (propeller disk loading) -- unit/description typed before
# the "DL" label itself.
$wsAircraft.Range("C19").Value = "N/m2"
$wsAircraft.Range("D19").Value = "Selected propeller disk loading"

# Row 20: Aspect ratio -- typed before going back to finish row 19.
$wsAircraft.Range("A20").Value = "A"
$wsAircraft.Range("D20").Value = "Aspect ratio"

# Back to row 19's "DL" label.
$wsAircraft.Range("A19").Value = "DL"

# Remaining cells: numbers, and "kg"/"-" strings already present in the
# shared-string table.
$wsAircraft.Range("B16").Value = 3353.95
$wsAircraft.Range("C16").Value = "kg"
$wsAircraft.Rows.Item(16).RowHeight = 23.25

$wsAircraft.Range("B17").Value = 4
$wsAircraft.Range("B17").NumberFormat = "0.000"
$wsAircraft.Range("C17").Value = "-"
$wsAircraft.Rows.Item(17).RowHeight = 23.25

$wsAircraft.Range("B18").Value = 2
$wsAircraft.Range("B18").NumberFormat = "0.000"
$wsAircraft.Range("C18").Value = "-"
$wsAircraft.Rows.Item(18).RowHeight = 23.25

$wsAircraft.Range("B19").Value = 1300
$wsAircraft.Range("B19").NumberFormat = "0.000"
$wsAircraft.Rows.Item(19).RowHeight = 23.25

$wsAircraft.Range("B20").Value = 5
$wsAircraft.Range("B20").NumberFormat = "0.000"
$wsAircraft.Range("C20").Value = "-"

# --- View state: make "Aircraft parameters" the active/front sheet ---
# (was "Mission"). "Mission" keeps its K8 selection but scrolls down one
# row (topLeftCell A2 -> A3); "Aircraft parameters" becomes the active
# tab, scrolled to A6, with the selection moved down to the newly added
# A20 (Aspect ratio) cell.
$wsMission.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1

$wsAircraft.Activate() | Out-Null
$wsAircraft.Range("A20").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
